# Fix header spacing when unchecked
#
# The default header's {{#SHOW_REPORT_UNDER}} block contains a manual
# line break (<w:br/>) between {{REPORT_UNDER_TITLE}} and
# {{REPORT_NUMBERS}}. When SHOW_REPORT_UNDER is false, the mail-merge
# tags collapse away but the forced line break remains, leaving a stray
# blank line in the header. Remove that manual line break.

$d = $word.ActiveDocument

# Locate the primary ("default") header for the first section - this is
# the one containing the SHOW_REPORT_UNDER merge tags.
$hdr = $d.Sections(1).Headers(1)
$rng = $hdr.Range

# "^l" matches a manual line break (w:br). Replace it with nothing so the
# run simply disappears, merging the surrounding text together.
$rng.Find.Execute("^l", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
